$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.141.01'
$ws.Range("E2").Value = '  -5.66%  '

$ws.Range("D3").Value = '3.348.57'
$ws.Range("E3").Value = '  -2.07%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '565.65'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.57%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '131.07'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.55%  '

$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").Value = '3.348.52'
$ws.Range("E8").Value = '  -2.06%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.473'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.66%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.43'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.74%  '

$ws.Range("E11").Value = '  -4.91%  '

$ws.Range("E12").Value = '  -1.44%  '

$ws.Range("D13").Value = '3.916.95'
$ws.Range("E13").Value = '  -2.12%  '

$ws.Range("E14").Value = '  -0.35%  '

$ws.Range("D15").Value = '3.340.64'
$ws.Range("E15").Value = '  -2.34%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000170'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.92%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '24.62'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.28%  '

$ws.Range("D18").Value = '60.188.61'

$ws.Range("E19").Value = '  +0.73%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.46'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.68%  '

$ws.Range("E21").Value = '  -7.45%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '354.48'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -7.59%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.559'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.85%  '

$ws.Range("D24").Value = '3.481.31'
$ws.Range("E24").Value = '  -2.17%  '

$ws.Range("E25").Value = '  +0.02%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '69.37'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -6.33%  '

$ws.Range("E27").Value = '  +2.13%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.66'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +17.75%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.53'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +7.34%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.998'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.14%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.94'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.61%  '

$ws.Range("E32").Value = '  +1.00%  '

$ws.Range("E33").Value = '  -3.23%  '

$ws.Range("E34").Value = '  -0.06%  '

$ws.Range("D35").Value = '3.379.19'
$ws.Range("E35").Value = '  -2.06%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '23.00'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.89%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.45'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.11%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.90'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.33%  '

$ws.Range("E39").Value = '  +0.11%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '158.96'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.05%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0768'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.37%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.14%  '

$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.38'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.48%  '

$ws.Range("B44").Value = 'ONDO'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.19'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +8.79%  '

$ws.Range("E45").Value = '  -4.50%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '40.77'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.56%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '24.02'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.53%  '

$ws.Range("E48").Value = '  -1.66%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.81'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.62%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.52'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +11.14%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.891'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.23%  '
